# Update to use 5V for consistent ATTINY20 programming
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ird-v17 bom")

# K1 - Relay: DC3 -> DC5 coil voltage part, price update
$ws.Range("C14").Value = "G6K-2G-Y-TR DC5"
$ws.Range("E14").Value = 4.21
$ws.Range("G14").Value = 39.79

# R2 - 27R -> 56R
$ws.Range("B17").Value = "56R"
$ws.Range("C17").Value = "RCS040256R0FKED"
$ws.Range("G17").Value = 0.57999999999999996

# R3 - 120R -> 750R
$ws.Range("B18").Value = "750R"
$ws.Range("C18").Value = "RCS0402750RFKED"
$ws.Range("E18").Value = 0.12
$ws.Range("G18").Value = 0.57999999999999996

# U2 - IR -> IR Detector (label clarification)
$ws.Range("B20").Value = "IR Detector"

# U3 - 3.3V LDO -> 5V SMPS module, consistent with 5V ATTINY20 programming
$ws.Range("B21").Value = "5V SMPS"
$ws.Range("C21").Value = "TLVM365R1RDNR"
$ws.Range("E21").Value = 2.2799999999999998
$ws.Range("G21").Value = 16.87

# VR1 - trimmer SKU correction
$ws.Range("C22").Value = "TC42X-2-103E"

# Move the active selection like the author's last saved state
$ws.Range("G16").Select()
